# "New domain name set up" - refresh currency rates report:
#  - add new currency rows (Euro, Franco Suiço, Iene, Libra, Peso Argentino,
#    Peso Chileno, Peso Colombiano, Peso Mexicano, Yuan)
#  - update existing USD/AUD/CAD values
#  - move the "Data das cotações" / "Horário do relatório" labels (and their
#    special formatting) down from row 6/7 to row 15/16, with fresh values
#  - restore the vacated C6/D6 cells to the ordinary body style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relocate the special "Data das cotações" / "Horário do relatório"
#     header formatting from row 6 down to row 15 before anything else
#     touches those cells. ---
$ws.Range("C6").Copy()
$ws.Range("C15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D6").Copy()
$ws.Range("D15").PasteSpecial(-4122)   # xlPasteFormats

# Restore C6/D6 back to the plain body-cell style used throughout the table.
$ws.Range("A6").Copy()
$ws.Range("C6").PasteSpecial(-4122)    # xlPasteFormats
$ws.Range("D6").PasteSpecial(-4122)    # xlPasteFormats

$excel.CutCopyMode = $false

# --- Update the existing three currency rows ---
$ws.Range("B2").Value = 4.6
$ws.Range("B3").Value = 3.47
$ws.Range("B4").Value = 3.68

# --- Row 5: Euro ---
$ws.Range("A5").Value = "Euro"
$ws.Range("B5").Value = 5.04
$ws.Range("C5").Value = "EUR"
$ws.Range("D5").Value = [char]0x20AC   # €

# --- Row 6: Franco Suíço ---
$ws.Range("A6").Value = "Franco Suíço"
$ws.Range("B6").Value = 4.96
$ws.Range("C6").Value = "CHF"
$ws.Range("D6").Value = "Fr"

# --- Row 7: Iene ---
$ws.Range("A7").Value = "Iene"
$ws.Range("B7").Value = 0.0374
$ws.Range("C7").Value = "JPY"
$ws.Range("D7").Value = [char]0x00A5   # ¥

# --- Row 8: Libra ---
$ws.Range("A8").Value = "Libra"
$ws.Range("B8").Value = 6.03
$ws.Range("C8").Value = "GBP"
$ws.Range("D8").Value = [char]0x00A3   # £

# --- Row 9: Peso Argentino ---
$ws.Range("A9").Value = "Peso Argentino"
$ws.Range("B9").Value = 0.041
$ws.Range("C9").Value = "ARS"
$ws.Range("D9").Value = "$"

# --- Row 10: Peso Chileno ---
$ws.Range("A10").Value = "Peso Chileno"
$ws.Range("B10").Value = 0.0059
$ws.Range("C10").Value = "CLP"
$ws.Range("D10").Value = "$"

# --- Row 11: Peso Colombiano ---
$ws.Range("A11").Value = "Peso Colombiano"
$ws.Range("B11").Value = 0.0012
$ws.Range("C11").Value = "COP"
$ws.Range("D11").Value = "$"

# --- Row 12: Peso Mexicano ---
$ws.Range("A12").Value = "Peso Mexicano"
$ws.Range("B12").Value = 0.23
$ws.Range("C12").Value = "MXN"
$ws.Range("D12").Value = "$"

# --- Row 13: Yuan ---
$ws.Range("A13").Value = "Yuan"
$ws.Range("B13").Value = 0.72
$ws.Range("C13").Value = "CNY"
$ws.Range("D13").Value = [char]0x00A5   # ¥

# --- Row 15: relocated header labels (formatting already copied above) ---
$ws.Range("C15").Value = "Data das cotações"
$ws.Range("D15").Value = "Horário do relatório"

# --- Row 16: refreshed quote date / report time ---
# A leading apostrophe stops "04/04/2022" being auto-parsed into a date
# serial; re-pasting the plain body format afterwards strips the
# resulting quote-prefix marker so the cell matches the rest of the table.
$ws.Range("C16").Value = "'04/04/2022"
$ws.Range("A16").Copy()
$ws.Range("C16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D16").Value = "20:35"
